$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataPackages")

# Bump the displayed format/version string (row 2, "Format:" / value cell)
# from "v0.0.1" to "v0.0.2".
$ws.Range("C2").Value = "v0.0.2"

# Freeze panes so that columns A:B and rows 1:7 stay fixed, with the
# view scrolled so C8 is the first cell of the scrollable area.
$ws.Activate()
$ws.Range("C8").Select()
$excel.ActiveWindow.FreezePanes = $true
